$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Reference cell that already carries the "text" number format (numFmtId 49,
# right-aligned) used throughout column C for Sucursal codes like "001".
$txtFmt = $ws.Range("C33").NumberFormat
$txtAlign = $ws.Range("C33").HorizontalAlignment

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = $txtFmt
    $cell.HorizontalAlignment = $txtAlign
    $cell.Value = $value
}

# Row 34
$ws.Range("A34").Value = "CCUENCA"
Set-TextCell $ws.Range("C34") "001"

# Row 35
$ws.Range("A35").Value = "F04033"
$ws.Range("C35").Value = 533

# Row 36
$ws.Range("A36").Value = "F04169"
Set-TextCell $ws.Range("C36") "369"

# Row 37
$ws.Range("A37").Value = "F00463"
Set-TextCell $ws.Range("C37") "063"

# Row 38
$ws.Range("A38").Value = "JANDINO"
$ws.Range("C38").Value = 102

# Row 39
$ws.Range("A39").Value = "F00219"
Set-TextCell $ws.Range("C39") "019"

# Row 40
$ws.Range("A40").Value = "F00089"
Set-TextCell $ws.Range("C40") "089"

# Make "Users" the active sheet/tab and park the selection where the
# author left it after adding the new rows.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("C41").Select()
